$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2374.75
$ws.Range("I40").Value = 1749.8334
$ws.Range("J40").Value = 2999.6667
$ws.Range("K40").Value = 1749.8334
$ws.Range("L40").Value = 2999.6667
$ws.Range("M40").Value = -1574.8334
$ws.Range("N40").Value = -3349.6667

$ws.Range("H76").Value = 6516.1177
$ws.Range("I76").Value = 5497
$ws.Range("J76").Value = 7229.5
$ws.Range("K76").Value = 5497
$ws.Range("L76").Value = 7229.5
$ws.Range("M76").Value = -5182

$ws.Range("H79").Value = 6516.1177
$ws.Range("I79").Value = 5497
$ws.Range("J79").Value = 7229.5
$ws.Range("K79").Value = 5497
$ws.Range("L79").Value = 7229.5
$ws.Range("M79").Value = -4405

$ws.Range("H88").Value = 1447.5454
$ws.Range("I88").Value = 394.25
$ws.Range("J88").Value = 2049.4285
$ws.Range("K88").Value = 394.25
$ws.Range("L88").Value = 2049.4285
$ws.Range("M88").Value = 11.75

$ws.Range("H91").Value = 1447.5454
$ws.Range("I91").Value = 394.25
$ws.Range("J91").Value = 2049.4285
$ws.Range("K91").Value = 394.25
$ws.Range("L91").Value = 2049.4285
$ws.Range("M91").Value = 1009.75

$ws.Range("H96").Value = 2005.2858
$ws.Range("I96").Value = 1207.2
$ws.Range("J96").Value = 4000.5
$ws.Range("K96").Value = 3621.6
$ws.Range("L96").Value = 12001.5
$ws.Range("M96").Value = -2248.6
$ws.Range("N96").Value = -14747.5

$ws.Range("H111").Value = 937.3333
$ws.Range("I111").Value = 397.33334
$ws.Range("J111").Value = 1477.3334
$ws.Range("K111").Value = 1192.00002
$ws.Range("L111").Value = 4432.0002
$ws.Range("M111").Value = 1874.99998

$ws.Range("H132").Value = 1501.875
$ws.Range("I132").Value = 1382.65
$ws.Range("J132").Value = 2098
$ws.Range("K132").Value = 4147.950000000001
$ws.Range("L132").Value = 6294
$ws.Range("M132").Value = -1617.950000000001
$ws.Range("N132").Value = -11354

$ws.Range("H135").Value = 1050.6111
$ws.Range("I135").Value = 540.3077
$ws.Range("J135").Value = 2377.4
$ws.Range("K135").Value = 4862.7693
$ws.Range("L135").Value = 21396.6
$ws.Range("M135").Value = -2327.7693
$ws.Range("N135").Value = -26466.6

$ws.Range("H137").Value = 2274
$ws.Range("I137").Value = 1240.7858
$ws.Range("J137").Value = 3589
$ws.Range("K137").Value = 3722.3574
$ws.Range("L137").Value = 10767
$ws.Range("M137").Value = -1172.3574
$ws.Range("N137").Value = -15867

$ws.Range("H138").Value = 3816.697
$ws.Range("I138").Value = 1162
$ws.Range("J138").Value = 4235.86
$ws.Range("K138").Value = 3486
$ws.Range("L138").Value = 12707.58
$ws.Range("M138").Value = 1654

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 639.1667
$ws.Range("I2").Value = 639.1667
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 639.1667
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -526.1667

$ws.Range("H61").Value = 1394.2858
$ws.Range("I61").Value = 1394.2858
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1394.2858
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1182.2858

$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("N97").ClearContents()

$ws.Range("H116").Value = 639.1667
$ws.Range("I116").Value = 639.1667
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 639.1667
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1654.8333

$ws.Range("H132").Value = 1094.3334
$ws.Range("I132").Value = 1094.3334
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3283.0002
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -753.0001999999999

$ws.Range("H136").Value = 1394.2858
$ws.Range("I136").Value = 1394.2858
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4182.857400000001
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -1632.857400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 639.1667
$ws.Range("I3").Value = 639.1667
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 639.1667
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -525.1667

$ws.Range("H94").Value = 899.8
$ws.Range("I94").Value = 624.75
$ws.Range("J94").Value = 2000
$ws.Range("K94").Value = 624.75
$ws.Range("L94").Value = 2000
$ws.Range("M94").Value = -173.75

$ws.Range("H105").Value = 3775.0322
$ws.Range("I105").Value = 3112
$ws.Range("J105").Value = 4824.8335
$ws.Range("K105").Value = 3112
$ws.Range("L105").Value = 4824.8335
$ws.Range("M105").Value = -1365

$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws.Range("H134").Value = 2889.3076
$ws.Range("I134").Value = 2614.2173
$ws.Range("J134").Value = 4998.3335
$ws.Range("K134").Value = 7842.651899999999
$ws.Range("L134").Value = 14995.0005
$ws.Range("M134").Value = -5307.651899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3314.4736
$ws.Range("I58").Value = 1902.2222
$ws.Range("J58").Value = 4585.5
$ws.Range("K58").Value = 1902.2222
$ws.Range("L58").Value = 4585.5
$ws.Range("M58").Value = -1699.2222
$ws.Range("N58").Value = -4991.5

$ws.Range("H59").Value = 43999
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 43999
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 43999
$ws.Range("N59").Value = -46289

$ws.Range("H86").Value = 13563.143
$ws.Range("I86").Value = 9648.333000000001
$ws.Range("J86").Value = 16499.25
$ws.Range("K86").Value = 9648.333000000001
$ws.Range("L86").Value = 16499.25
$ws.Range("M86").Value = -8525.333000000001

$ws.Range("H89").Value = 13563.143
$ws.Range("I89").Value = 9648.333000000001
$ws.Range("J89").Value = 16499.25
$ws.Range("K89").Value = 48241.665
$ws.Range("L89").Value = 82496.25
$ws.Range("M89").Value = -42625.665

$ws.Range("H132").Value = 1866.3572
$ws.Range("I132").Value = 1560.5333
$ws.Range("J132").Value = 2219.2307
$ws.Range("K132").Value = 4681.5999
$ws.Range("L132").Value = 6657.6921
$ws.Range("M132").Value = -2151.5999
$ws.Range("N132").Value = -11717.6921

$ws.Range("H134").Value = 3956.5293
$ws.Range("I134").Value = 3391.8333
$ws.Range("J134").Value = 5311.8
$ws.Range("K134").Value = 10175.4999
$ws.Range("L134").Value = 15935.4
$ws.Range("M134").Value = -7640.499899999999

$ws.Range("H136").Value = 3314.4736
$ws.Range("I136").Value = 1902.2222
$ws.Range("J136").Value = 4585.5
$ws.Range("K136").Value = 5706.6666
$ws.Range("L136").Value = 13756.5
$ws.Range("M136").Value = -3156.6666
$ws.Range("N136").Value = -18856.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 657.3
$ws.Range("I2").Value = 89.666664
$ws.Range("J2").Value = 1508.75
$ws.Range("K2").Value = 89.666664
$ws.Range("L2").Value = 1508.75
$ws.Range("M2").Value = 23.333336
$ws.Range("N2").Value = -1734.75

$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()

$ws.Range("H80").Value = 2972.6
$ws.Range("I80").Value = 2165.5
$ws.Range("J80").Value = 3510.6667
$ws.Range("K80").Value = 2165.5
$ws.Range("L80").Value = 3510.6667
$ws.Range("M80").Value = -1167.5

$ws.Range("H83").Value = 2972.6
$ws.Range("I83").Value = 2165.5
$ws.Range("J83").Value = 3510.6667
$ws.Range("K83").Value = 10827.5
$ws.Range("L83").Value = 17553.3335
$ws.Range("M83").Value = -5835.5

$ws.Range("H102").Value = 1279.0454
$ws.Range("I102").Value = 452.16666
$ws.Range("J102").Value = 5000
$ws.Range("K102").Value = 452.16666
$ws.Range("L102").Value = 5000
$ws.Range("M102").Value = 1169.83334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 86.40000000000001
$ws.Range("I93").Value = 107
$ws.Range("J93").Value = 4
$ws.Range("K93").Value = 107
$ws.Range("L93").Value = 4
$ws.Range("M93").Value = 1141
$ws.Range("N93").Value = -2500

$ws.Range("H122").Value = 6541.8335
$ws.Range("I122").Value = 3055.7778
$ws.Range("J122").Value = 17000
$ws.Range("K122").Value = 9167.3334
$ws.Range("L122").Value = 51000
$ws.Range("M122").Value = -6717.3334
$ws.Range("N122").Value = -55900

$ws.Range("H132").Value = 4769.8975
$ws.Range("I132").Value = 4263.25
$ws.Range("J132").Value = 5580.533
$ws.Range("K132").Value = 12789.75
$ws.Range("L132").Value = 16741.599
$ws.Range("M132").Value = -10259.75
$ws.Range("N132").Value = -21801.599

$ws.Range("H136").Value = 2356.6155
$ws.Range("I136").Value = 2511.3333
$ws.Range("J136").Value = 500
$ws.Range("K136").Value = 7533.999899999999
$ws.Range("L136").Value = 1500
$ws.Range("M136").Value = -4983.999899999999
$ws.Range("N136").Value = -6600

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1078.7333
$ws.Range("I81").Value = 1109.4
$ws.Range("J81").Value = 1017.4
$ws.Range("K81").Value = 2218.8
$ws.Range("L81").Value = 2034.8
$ws.Range("M81").Value = -1157.8

$ws.Range("H84").Value = 1078.7333
$ws.Range("I84").Value = 1109.4
$ws.Range("J84").Value = 1017.4
$ws.Range("K84").Value = 11094
$ws.Range("L84").Value = 10174
$ws.Range("M84").Value = -5790

$ws.Range("H107").Value = 1040.5
$ws.Range("I107").Value = 435.875
$ws.Range("J107").Value = 2249.75
$ws.Range("K107").Value = 1307.625
$ws.Range("L107").Value = 6749.25
$ws.Range("M107").Value = 612.375

$ws.Range("H126").Value = 145632.58
$ws.Range("I126").Value = 200485.8
$ws.Range("J126").Value = 8499.5
$ws.Range("K126").Value = 601457.3999999999
$ws.Range("L126").Value = 25498.5
$ws.Range("M126").Value = -598987.3999999999

$ws.Range("H132").Value = 1153.1578
$ws.Range("I132").Value = 1097.6428
$ws.Range("J132").Value = 1308.6
$ws.Range("K132").Value = 3292.9284
$ws.Range("L132").Value = 3925.8
$ws.Range("M132").Value = -762.9284000000002

$ws.Range("H136").Value = 60588.293
$ws.Range("I136").Value = 1026.8462
$ws.Range("J136").Value = 254163
$ws.Range("K136").Value = 3080.5386
$ws.Range("L136").Value = 762489
$ws.Range("M136").Value = -530.5385999999999
